# Tweak to title slide for bi-implication
# Slide 1, Shape 3 ("Title 1", the ctrTitle placeholder) currently holds a
# single run "Bi-implication" at sz=6600. Split it into three runs so the
# "<->" arrow in the middle renders smaller than the surrounding text:
#   "Bi"            sz=6600
#   "\u2194" (<->)  sz=4000
#   "implication"   sz=6600

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tr  = $shp.TextFrame.TextRange

$arrow = [char]0x2194
$tr.Text = "Bi" + $arrow + "implication"

# "Bi" -> characters 1-2 (inherits the sz=6600 already on the run; set
# explicitly so it is an independent run with the right size)
$runBi = $tr.Characters(1, 2)
$runBi.Font.Size = 66

# the bi-implication arrow -> character 3, shown smaller
$runArrow = $tr.Characters(3, 1)
$runArrow.Font.Size = 40

# "implication" -> characters 4-14
$runImpl = $tr.Characters(4, 11)
$runImpl.Font.Size = 66
